$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for rows 2-32.
# It moves from 45183 (2023-09-14) to 45184 (2023-09-15) for every row.
$ws.Range("C2:C32").Value = 45184
